$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:06"
$ws.Cells.Item(3,1).Value = "Total filas: 447"
$ws.Cells.Item(61,1).Value = "06:02:16"
$ws.Cells.Item(61,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(61,4).Value = 90
$ws.Cells.Item(62,1).Value = "05:47:32"
$ws.Cells.Item(62,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(62,4).Value = 105
$ws.Cells.Item(214,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(215,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(229,3).Value = "16_SANTA ANA"
$ws.Cells.Item(230,3).Value = "17_ROMERO"
$ws.Cells.Item(269,1).Value = "13:24:27"
$ws.Cells.Item(269,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(269,4).Value = 40
$ws.Cells.Item(270,1).Value = "12:44:05"
$ws.Cells.Item(270,3).Value = "17_ROMERO"
$ws.Cells.Item(270,4).Value = 80
$ws.Cells.Item(271,1).Value = "13:51:56"
$ws.Cells.Item(271,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(271,4).Value = 13
$ws.Cells.Item(298,1).Value = "14:56:04"
$ws.Cells.Item(298,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(298,4).Value = 9
$ws.Cells.Item(299,1).Value = "13:24:27"
$ws.Cells.Item(299,3).Value = "10_OLMOS"
$ws.Cells.Item(299,4).Value = 101
$ws.Cells.Item(326,3).Value = "17_ROMERO"
$ws.Cells.Item(327,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(338,1).Value = "15:53:26"
$ws.Cells.Item(338,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(338,4).Value = 12
$ws.Cells.Item(339,1).Value = "14:42:26"
$ws.Cells.Item(339,3).Value = "27_EL RETIRO"
$ws.Cells.Item(339,4).Value = 83
$ws.Cells.Item(421,1).Value = "17:48:06"
$ws.Cells.Item(421,2).Value = "17:53"
$ws.Cells.Item(421,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(421,4).Value = 5
$ws.Cells.Item(422,1).Value = "17:48:06"
$ws.Cells.Item(422,2).Value = "18:01"
$ws.Cells.Item(422,3).Value = "16_SANTA ANA"
$ws.Cells.Item(422,4).Value = 13
$ws.Cells.Item(423,1).Value = "17:34:37"
$ws.Cells.Item(423,2).Value = "18:03"
$ws.Cells.Item(423,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(423,4).Value = 29
$ws.Cells.Item(424,1).Value = "16:30:20"
$ws.Cells.Item(424,2).Value = "18:04"
$ws.Cells.Item(424,3).Value = "17_ROMERO"
$ws.Cells.Item(424,4).Value = 94
$ws.Cells.Item(425,1).Value = "17:48:06"
$ws.Cells.Item(425,2).Value = "18:05"
$ws.Cells.Item(425,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(425,4).Value = 17
$ws.Cells.Item(426,1).Value = "17:48:06"
$ws.Cells.Item(426,2).Value = "18:11"
$ws.Cells.Item(426,3).Value = "16_SANTA ANA"
$ws.Cells.Item(426,4).Value = 23
$ws.Cells.Item(427,1).Value = "17:48:06"
$ws.Cells.Item(427,2).Value = "18:16"
$ws.Cells.Item(427,3).Value = "15_ABASTO"
$ws.Cells.Item(427,4).Value = 28
$ws.Cells.Item(428,1).Value = "17:48:06"
$ws.Cells.Item(428,2).Value = "18:17"
$ws.Cells.Item(428,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(428,4).Value = 29
$ws.Cells.Item(429,1).Value = "16:30:20"
$ws.Cells.Item(429,2).Value = "18:21"
$ws.Cells.Item(429,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(429,4).Value = 111
$ws.Cells.Item(430,1).Value = "16:53:43"
$ws.Cells.Item(430,2).Value = "18:24"
$ws.Cells.Item(430,3).Value = "14_ABASTO"
$ws.Cells.Item(430,4).Value = 91
$ws.Cells.Item(431,1).Value = "17:48:06"
$ws.Cells.Item(431,2).Value = "18:25"
$ws.Cells.Item(431,3).Value = "14_ABASTO"
$ws.Cells.Item(431,4).Value = 37
$ws.Cells.Item(432,1).Value = "16:39:47"
$ws.Cells.Item(432,2).Value = "18:27"
$ws.Cells.Item(432,3).Value = "215C_EL PATO"
$ws.Cells.Item(432,4).Value = 108
$ws.Cells.Item(433,1).Value = "16:30:20"
$ws.Cells.Item(433,2).Value = "18:28"
$ws.Cells.Item(433,3).Value = "215C_EL PATO"
$ws.Cells.Item(433,4).Value = 118
$ws.Cells.Item(434,1).Value = "16:39:47"
$ws.Cells.Item(434,2).Value = "18:32"
$ws.Cells.Item(434,3).Value = "11X44_ETCHEVERRY"
$ws.Cells.Item(434,4).Value = 113
$ws.Cells.Item(435,1).Value = "17:48:06"
$ws.Cells.Item(435,2).Value = "18:40"
$ws.Cells.Item(435,3).Value = "15_ABASTO"
$ws.Cells.Item(435,4).Value = 52
$ws.Cells.Item(436,2).Value = "18:45"
$ws.Cells.Item(436,3).Value = "14_ABASTO"
$ws.Cells.Item(436,4).Value = 71
$ws.Cells.Item(437,1).Value = "16:53:43"
$ws.Cells.Item(437,2).Value = "18:48"
$ws.Cells.Item(437,3).Value = "14X44_ABASTO"
$ws.Cells.Item(437,4).Value = 115
$ws.Cells.Item(438,1).Value = "17:34:37"
$ws.Cells.Item(438,2).Value = "18:56"
$ws.Cells.Item(438,3).Value = "10_OLMOS"
$ws.Cells.Item(438,4).Value = 82
$ws.Cells.Item(438,5).Value = "LP1912"
$ws.Cells.Item(439,1).Value = "17:13:46"
$ws.Cells.Item(439,2).Value = "18:58"
$ws.Cells.Item(439,3).Value = "215A_EL PATO"
$ws.Cells.Item(439,4).Value = 105
$ws.Cells.Item(439,5).Value = "LP1912"
$ws.Cells.Item(440,1).Value = "17:48:06"
$ws.Cells.Item(440,2).Value = "18:59"
$ws.Cells.Item(440,3).Value = "215A_EL PATO"
$ws.Cells.Item(440,4).Value = 71
$ws.Cells.Item(440,5).Value = "LP1912"
$ws.Cells.Item(441,1).Value = "17:13:46"
$ws.Cells.Item(441,2).Value = "19:04"
$ws.Cells.Item(441,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(441,4).Value = 111
$ws.Cells.Item(441,5).Value = "LP1912"
$ws.Cells.Item(442,1).Value = "17:48:06"
$ws.Cells.Item(442,2).Value = "19:05"
$ws.Cells.Item(442,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(442,4).Value = 77
$ws.Cells.Item(442,5).Value = "LP1912"
$ws.Cells.Item(443,1).Value = "17:13:46"
$ws.Cells.Item(443,2).Value = "19:05"
$ws.Cells.Item(443,3).Value = "27_EL RETIRO"
$ws.Cells.Item(443,4).Value = 112
$ws.Cells.Item(443,5).Value = "LP1912"
$ws.Cells.Item(444,1).Value = "17:13:46"
$ws.Cells.Item(444,2).Value = "19:10"
$ws.Cells.Item(444,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(444,4).Value = 117
$ws.Cells.Item(444,5).Value = "LP1912"
$ws.Cells.Item(445,1).Value = "17:48:06"
$ws.Cells.Item(445,2).Value = "19:11"
$ws.Cells.Item(445,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(445,4).Value = 83
$ws.Cells.Item(445,5).Value = "LP1912"
$ws.Cells.Item(446,1).Value = "17:48:06"
$ws.Cells.Item(446,2).Value = "19:12"
$ws.Cells.Item(446,3).Value = "10_OLMOS"
$ws.Cells.Item(446,4).Value = 84
$ws.Cells.Item(446,5).Value = "LP1912"
$ws.Cells.Item(447,1).Value = "17:34:37"
$ws.Cells.Item(447,2).Value = "19:16"
$ws.Cells.Item(447,3).Value = "17_ROMERO"
$ws.Cells.Item(447,4).Value = 102
$ws.Cells.Item(447,5).Value = "LP1912"
$ws.Cells.Item(448,1).Value = "17:34:37"
$ws.Cells.Item(448,2).Value = "19:16"
$ws.Cells.Item(448,3).Value = "27_EL RETIRO"
$ws.Cells.Item(448,4).Value = 102
$ws.Cells.Item(448,5).Value = "LP1912"
$ws.Cells.Item(449,1).Value = "17:48:06"
$ws.Cells.Item(449,2).Value = "19:17"
$ws.Cells.Item(449,3).Value = "27_EL RETIRO"
$ws.Cells.Item(449,4).Value = 89
$ws.Cells.Item(449,5).Value = "LP1912"
$ws.Cells.Item(450,1).Value = "17:48:06"
$ws.Cells.Item(450,2).Value = "19:19"
$ws.Cells.Item(450,3).Value = "17_ROMERO"
$ws.Cells.Item(450,4).Value = 91
$ws.Cells.Item(450,5).Value = "LP1912"
$ws.Cells.Item(451,1).Value = "17:34:37"
$ws.Cells.Item(451,2).Value = "19:30"
$ws.Cells.Item(451,3).Value = "225_GOMEZ"
$ws.Cells.Item(451,4).Value = 116
$ws.Cells.Item(451,5).Value = "LP1912"
$ws.Cells.Item(452,1).Value = "17:48:06"
$ws.Cells.Item(452,2).Value = "19:40"
$ws.Cells.Item(452,3).Value = "215C_EL PATO"
$ws.Cells.Item(452,4).Value = 112
$ws.Cells.Item(452,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:06"
$ws.Cells.Item(3,1).Value = "Total filas: 44"
$ws.Cells.Item(48,1).Value = "17:48:06"
$ws.Cells.Item(48,2).Value = "18:59"
$ws.Cells.Item(48,3).Value = "215A_EL PATO"
$ws.Cells.Item(48,4).Value = 71
$ws.Cells.Item(48,5).Value = "LP1912"
$ws.Cells.Item(49,1).Value = "17:48:06"
$ws.Cells.Item(49,2).Value = "19:40"
$ws.Cells.Item(49,3).Value = "215C_EL PATO"
$ws.Cells.Item(49,4).Value = 112
$ws.Cells.Item(49,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 17:48:06"
$ws.Cells.Item(3,1).Value = "Total filas: 53"
$ws.Cells.Item(58,1).Value = "17:48:06"
$ws.Cells.Item(58,2).Value = "19:05"
$ws.Cells.Item(58,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(58,4).Value = 77
$ws.Cells.Item(58,5).Value = "L6173"
